$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.049481999999999
$ws.Range("H2").Value = 27.148446
$ws.Range("I2").Value = 0.2715881048104983
$ws.Range("J2").Value = 0.2887858053066977
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.024115666666666
$ws.Range("N2").Value = 6.072347
$ws.Range("O2").Value = 0.04033801028669426
$ws.Range("P2").Value = 0.05315356654931184
$ws.Range("Q2").Value = 18.31719829141799
$ws.Range("R2").Value = 164.854784622762
$ws.Range("S2").Value = 0.01095532376558968
$ws.Range("T2").Value = 0.01534999552086617

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.049481999999999
$ws.Range("H3").Value = 27.148446
$ws.Range("I3").Value = 0.2715881048104983
$ws.Range("J3").Value = 0.2887858053066977
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.85972166666667
$ws.Range("N3").Value = 35.579165
$ws.Range("O3").Value = 0.2363489312718777
$ws.Range("P3").Value = 0.3114379851145606
$ws.Range("Q3").Value = 107.32433774751
$ws.Range("R3").Value = 965.9190397275901
$ws.Range("S3").Value = 0.06418955831811599
$ws.Range("T3").Value = 0.08993886933440373

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.049481999999999
$ws.Range("H4").Value = 27.148446
$ws.Range("I4").Value = 0.2715881048104983
$ws.Range("J4").Value = 0.2887858053066977
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 36.29503
$ws.Range("N4").Value = 72.59005999999999
$ws.Range("O4").Value = 0.7233130584414281
$ws.Range("P4").Value = 0.6354084483361276
$ws.Range("Q4").Value = 328.4512206744599
$ws.Range("R4").Value = 1970.70732404676
$ws.Range("S4").Value = 0.1964432227267927
$ws.Range("T4").Value = 0.1834969404514279

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.313100666666667
$ws.Range("H5").Value = 15.939302
$ws.Range("I5").Value = 0.1594538715837432
$ws.Range("J5").Value = 0.1695509261965366
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.024115666666666
$ws.Range("N5").Value = 6.072347
$ws.Range("O5").Value = 0.04033801028669426
$ws.Range("P5").Value = 0.05315356654931184
$ws.Range("Q5").Value = 10.75433029797711
$ws.Range("R5").Value = 96.78897268179401
$ws.Range("S5").Value = 0.006432051912198258
$ws.Range("T5").Value = 0.009012236439085065

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.313100666666667
$ws.Range("H6").Value = 15.939302
$ws.Range("I6").Value = 0.1594538715837432
$ws.Range("J6").Value = 0.1695509261965366
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 11.85972166666667
$ws.Range("N6").Value = 35.579165
$ws.Range("O6").Value = 0.2363489312718777
$ws.Range("P6").Value = 0.3114379851145606
$ws.Range("Q6").Value = 63.01189509364779
$ws.Range("R6").Value = 567.1070558428301
$ws.Range("S6").Value = 0.03768675213598092
$ws.Range("T6").Value = 0.05280459882895692

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.313100666666667
$ws.Range("H7").Value = 15.939302
$ws.Range("I7").Value = 0.1594538715837432
$ws.Range("J7").Value = 0.1695509261965366
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 36.29503
$ws.Range("N7").Value = 72.59005999999999
$ws.Range("O7").Value = 0.7233130584414281
$ws.Range("P7").Value = 0.6354084483361276
$ws.Range("Q7").Value = 192.8391480896867
$ws.Range("R7").Value = 1157.03488853812
$ws.Range("S7").Value = 0.115335067535564
$ws.Range("T7").Value = 0.1077340909284946

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.671725333333334
$ws.Range("H8").Value = 14.015176
$ws.Range("I8").Value = 0.1402052658345742
$ws.Range("J8").Value = 0.1490834461639205
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.024115666666666
$ws.Range("N8").Value = 6.072347
$ws.Range("O8").Value = 0.04033801028669426
$ws.Range("P8").Value = 0.05315356654931184
$ws.Range("Q8").Value = 9.456112437563554
$ws.Range("R8").Value = 85.105011938072
$ws.Range("S8").Value = 0.005655601455483755
$ws.Range("T8").Value = 0.007924316877074697

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.671725333333334
$ws.Range("H9").Value = 14.015176
$ws.Range("I9").Value = 0.1402052658345742
$ws.Range("J9").Value = 0.1490834461639205
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 11.85972166666667
$ws.Range("N9").Value = 35.579165
$ws.Range("O9").Value = 0.2363489312718777
$ws.Range("P9").Value = 0.3114379851145606
$ws.Range("Q9").Value = 55.40536215644889
$ws.Range("R9").Value = 498.6482594080401
$ws.Range("S9").Value = 0.0331373647386911
$ws.Range("T9").Value = 0.04643024808722647

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.671725333333334
$ws.Range("H10").Value = 14.015176
$ws.Range("I10").Value = 0.1402052658345742
$ws.Range("J10").Value = 0.1490834461639205
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 36.29503
$ws.Range("N10").Value = 72.59005999999999
$ws.Range("O10").Value = 0.7233130584414281
$ws.Range("P10").Value = 0.6354084483361276
$ws.Range("Q10").Value = 169.5604111250933
$ws.Range("R10").Value = 1017.36246675056
$ws.Range("S10").Value = 0.1014122996403993
$ws.Range("T10").Value = 0.09472888119961935

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 8.333402
$ws.Range("H11").Value = 25.000206
$ws.Range("I11").Value = 0.2500975034597578
$ws.Range("J11").Value = 0.2659343603881907
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.024115666666666
$ws.Range("N11").Value = 6.072347
$ws.Range("O11").Value = 0.04033801028669426
$ws.Range("P11").Value = 0.05315356654931184
$ws.Range("Q11").Value = 16.86776954483133
$ws.Range("R11").Value = 151.809925903482
$ws.Range("S11").Value = 0.01008843566723626
$ws.Range("T11").Value = 0.01413535972264237

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 8.333402
$ws.Range("H12").Value = 25.000206
$ws.Range("I12").Value = 0.2500975034597578
$ws.Range("J12").Value = 0.2659343603881907
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 11.85972166666667
$ws.Range("N12").Value = 35.579165
$ws.Range("O12").Value = 0.2363489312718777
$ws.Range("P12").Value = 0.3114379851145606
$ws.Range("Q12").Value = 98.83182825644333
$ws.Range("R12").Value = 889.48645430799
$ws.Range("S12").Value = 0.05911027765647849
$ws.Range("T12").Value = 0.08282206137202755

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 8.333402
$ws.Range("H13").Value = 25.000206
$ws.Range("I13").Value = 0.2500975034597578
$ws.Range("J13").Value = 0.2659343603881907
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 36.29503
$ws.Range("N13").Value = 72.59005999999999
$ws.Range("O13").Value = 0.7233130584414281
$ws.Range("P13").Value = 0.6354084483361276
$ws.Range("Q13").Value = 302.46107559206
$ws.Range("R13").Value = 1814.76645355236
$ws.Range("S13").Value = 0.1808987901360431
$ws.Range("T13").Value = 0.1689769392935208

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 5.9529025
$ws.Range("H14").Value = 11.905805
$ws.Range("I14").Value = 0.1786552543114266
$ws.Range("J14").Value = 0.1266454619446545
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.024115666666666
$ws.Range("N14").Value = 6.072347
$ws.Range("O14").Value = 0.04033801028669426
$ws.Range("P14").Value = 0.05315356654931184
$ws.Range("Q14").Value = 12.04936321238917
$ws.Range("R14").Value = 72.296179274335
$ws.Range("S14").Value = 0.007206597486186306
$ws.Range("T14").Value = 0.006731657989643535

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 5.9529025
$ws.Range("H15").Value = 11.905805
$ws.Range("I15").Value = 0.1786552543114266
$ws.Range("J15").Value = 0.1266454619446545
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 11.85972166666667
$ws.Range("N15").Value = 35.579165
$ws.Range("O15").Value = 0.2363489312718777
$ws.Range("P15").Value = 0.3114379851145606
$ws.Range("Q15").Value = 70.59976675880418
$ws.Range("R15").Value = 423.5986005528251
$ws.Range("S15").Value = 0.0422249784226112
$ws.Range("T15").Value = 0.03944220749194597

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 5.9529025
$ws.Range("H16").Value = 11.905805
$ws.Range("I16").Value = 0.1786552543114266
$ws.Range("J16").Value = 0.1266454619446545
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 36.29503
$ws.Range("N16").Value = 72.59005999999999
$ws.Range("O16").Value = 0.7233130584414281
$ws.Range("P16").Value = 0.6354084483361276
$ws.Range("Q16").Value = 216.060774824575
$ws.Range("R16").Value = 864.2430992983
$ws.Range("S16").Value = 0.1292236784026291
$ws.Range("T16").Value = 0.08047159646306504

